# Apply the "dSF" (column F) updates described by the commit:
# "repull data, push all data, mean calculation"
#
# This updates specific cells in column F (dSF) on the active worksheet
# to reflect the repulled / recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 8
    3  = -7
    4  = 2
    5  = 0
    6  = 8
    7  = 10
    8  = -7
    10 = -6
    14 = 7
    17 = -10
    18 = -6
    20 = -2
    22 = -10
    23 = -7
    26 = -3
    27 = 7
    28 = -5
    30 = -1
    32 = -2
    34 = -2
    35 = -5
    36 = -5
    38 = -4
    41 = 2
    42 = -3
    43 = -7
    45 = -7
    47 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
